# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Vega Modelo de Temuco - Ají"
# at row 432 (pushing the existing rows 432..460 down to 433..461) and
# populate the new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 432..460 down to 433..461, leaving a blank row 432.
$ws.Rows.Item(432).Insert()

# Populate the newly inserted row 432 with the new record.
$ws.Cells.Item(432, 1).Value  = 10
$ws.Cells.Item(432, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(432, 3).Value  = "La Araucanía"
$ws.Cells.Item(432, 4).Value  = 44461
$ws.Cells.Item(432, 5).Value  = 9
$ws.Cells.Item(432, 6).Value  = 100112021
$ws.Cells.Item(432, 7).Value  = "Ají"
$ws.Cells.Item(432, 8).Value  = "Inferno"
$ws.Cells.Item(432, 9).Value  = "Extra"
$ws.Cells.Item(432, 10).Value = 20
$ws.Cells.Item(432, 11).Value = 55000
$ws.Cells.Item(432, 12).Value = 55000
$ws.Cells.Item(432, 13).Value = 55000
$ws.Cells.Item(432, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(432, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(432, 16).Value = 3667
$ws.Cells.Item(432, 17).Value = 15
$ws.Cells.Item(432, 18).Value = "Hortaliza"
